$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete columns F and G (product-sales-ratio / sales-volume point-in-time columns)
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(6).Delete()

# Step 2: for every year block of 4 quarters (A/B/C/D), swap the B-quarter and C-quarter rows
# (columns A-E) so that the C-quarter row appears before the B-quarter row in the sheet.
$ws.Range("A3").Value = "2000年C"
$ws.Range("A4").Value = "2000年B"
$ws.Range("B3").Value = 97.7
$ws.Range("B4").Value = 97.5
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D3").Value = 23.9
$ws.Range("D4").Value = 11.6
$ws.Range("E3").Value = 151.8
$ws.Range("E4").Value = 107.2
$ws.Range("A7").Value = "2001年C"
$ws.Range("A8").Value = "2001年B"
$ws.Range("B7").Value = 100.7
$ws.Range("B8").Value = 100.4
$ws.Range("C7").Value = 3
$ws.Range("C8").Value = ""
$ws.Range("D7").Value = 5.4
$ws.Range("D8").Value = -6.4
$ws.Range("E7").Value = 152.7
$ws.Range("E8").Value = 107.1
$ws.Range("A11").Value = "2002年C"
$ws.Range("A12").Value = "2002年B"
$ws.Range("B11").Value = 101.4
$ws.Range("B12").Value = 100.9
$ws.Range("C11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D11").Value = -11.4
$ws.Range("D12").Value = -7.3
$ws.Range("E11").Value = 149.6
$ws.Range("E12").Value = 105.6
$ws.Range("A15").Value = "2003年C"
$ws.Range("A16").Value = "2003年B"
$ws.Range("B15").Value = 99.5
$ws.Range("B16").Value = 100.5
$ws.Range("C15").Value = -1.9
$ws.Range("C16").Value = -0.4
$ws.Range("D15").Value = 5.2
$ws.Range("D16").Value = -9.1
$ws.Range("E15").Value = 137.9
$ws.Range("E16").Value = 118.1
$ws.Range("A19").Value = "2004年C"
$ws.Range("A20").Value = "2004年B"
$ws.Range("B19").Value = 102.8
$ws.Range("B20").Value = 97.5
$ws.Range("C19").Value = 3.3
$ws.Range("C20").Value = -3
$ws.Range("D19").Value = 19.3
$ws.Range("D20").Value = 1.2
$ws.Range("E19").Value = 130.7
$ws.Range("E20").Value = 84.59999999999999
$ws.Range("A23").Value = "2005年C"
$ws.Range("A24").Value = "2005年B"
$ws.Range("B23").Value = 99
$ws.Range("B24").Value = 96.5
$ws.Range("C23").Value = -3.8
$ws.Range("C24").Value = -1
$ws.Range("D23").Value = 4.5
$ws.Range("D24").Value = -2.2
$ws.Range("E23").Value = 144.8
$ws.Range("E24").Value = 94.09999999999999
$ws.Range("A27").Value = "2006年C"
$ws.Range("A28").Value = "2006年B"
$ws.Range("B27").Value = 99.59999999999999
$ws.Range("B28").Value = 99.90000000000001
$ws.Range("C27").Value = 0.6
$ws.Range("C28").Value = 3.4
$ws.Range("D27").Value = 9.1
$ws.Range("D28").Value = -10.4
$ws.Range("E27").Value = 150.5
$ws.Range("E28").Value = 101.1
$ws.Range("A31").Value = "2007年C"
$ws.Range("A32").Value = "2007年B"
$ws.Range("B31").Value = 96.8
$ws.Range("B32").Value = 98.7
$ws.Range("C31").Value = -2.8
$ws.Range("C32").Value = -1.2
$ws.Range("D31").Value = -6.5
$ws.Range("D32").Value = 12.9
$ws.Range("E31").Value = 151.3
$ws.Range("E32").Value = 96.7
$ws.Range("A35").Value = "2008年C"
$ws.Range("A36").Value = "2008年B"
$ws.Range("B35").Value = 96.3
$ws.Range("B36").Value = 93.09999999999999
$ws.Range("C35").Value = -3
$ws.Range("C36").Value = -4.9
$ws.Range("D35").Value = 81.3
$ws.Range("D36").Value = 129.7
$ws.Range("E35").Value = 147.2
$ws.Range("E36").Value = 114.3
$ws.Range("A39").Value = "2009年C"
$ws.Range("A40").Value = "2009年B"
$ws.Range("B39").Value = 100
$ws.Range("B40").Value = 101.5
$ws.Range("C39").Value = 1.1
$ws.Range("C40").Value = 2
$ws.Range("D39").Value = -2.6
$ws.Range("D40").Value = -10.5
$ws.Range("E39").Value = 141.3
$ws.Range("E40").Value = 91.5
$ws.Range("A43").Value = "2010年C"
$ws.Range("A44").Value = "2010年B"
$ws.Range("B43").Value = 98.5
$ws.Range("B44").Value = 95.3
$ws.Range("C43").Value = -1.2
$ws.Range("C44").Value = -3.9
$ws.Range("D43").Value = 23.8
$ws.Range("D44").Value = 31.6
$ws.Range("E43").Value = 169.7
$ws.Range("E44").Value = 114.7
$ws.Range("A47").Value = "2011年C"
$ws.Range("A48").Value = "2011年B"
$ws.Range("B47").Value = 99.2
$ws.Range("B48").Value = 97.40000000000001
$ws.Range("C47").Value = 0.8
$ws.Range("C48").Value = -1.1
$ws.Range("D47").Value = 9.800000000000001
$ws.Range("D48").Value = 39.8
$ws.Range("E47").Value = 183
$ws.Range("E48").Value = 127.8
$ws.Range("A51").Value = "2012年C"
$ws.Range("A52").Value = "2012年B"
$ws.Range("B51").Value = 98.2
$ws.Range("B52").Value = 97.59999999999999
$ws.Range("C51").Value = 5.8
$ws.Range("C52").Value = 0.5
$ws.Range("D51").Value = 30.01
$ws.Range("D52").Value = 30
$ws.Range("E51").Value = 127.4
$ws.Range("E52").Value = 84.59999999999999
$ws.Range("A55").Value = "2013年C"
$ws.Range("A56").Value = "2013年B"
$ws.Range("B55").Value = 100.3
$ws.Range("B56").Value = 97.90000000000001
$ws.Range("C55").Value = 2.4
$ws.Range("C56").Value = 0.8
$ws.Range("D55").Value = -5.3
$ws.Range("D56").Value = 25.2
$ws.Range("E55").Value = 143
$ws.Range("E56").Value = 95.3
$ws.Range("A59").Value = "2014年C"
$ws.Range("A60").Value = "2014年B"
$ws.Range("B59").Value = 99.3
$ws.Range("B60").Value = 99.7
$ws.Range("C59").Value = -0.3
$ws.Range("C60").Value = 2.4
$ws.Range("D59").Value = 9.5
$ws.Range("D60").Value = 3.7
$ws.Range("E59").Value = 125.5
$ws.Range("E60").Value = 86.5
$ws.Range("A63").Value = "2015年C"
$ws.Range("A64").Value = "2015年B"
$ws.Range("B63").Value = 99.3
$ws.Range("B64").Value = 100.6
$ws.Range("C63").Value = -0.2
$ws.Range("C64").Value = 0.5
$ws.Range("D63").Value = 9.4
$ws.Range("D64").Value = 2.4
$ws.Range("E63").Value = 102.1
$ws.Range("E64").Value = 70.59999999999999
